$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.380.21'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '1.795.95'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '307.65'
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("D7").Value = '0.4510'
$ws.Range("E7").Value = '  -1.53%  '
$ws.Range("E8").Value = '  -2.61%  '
$ws.Range("D9").Value = '46.07'
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '0.07079'
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '0.07801'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = '19.44'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '1.840.50'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '5.286'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").Value = '6.342'
$ws.Range("E16").Value = '  -0.88%  '
$ws.Range("D17").Value = '84.90'
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D19").Value = '0.000008521'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = '26.392.89'
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("D23").Value = '4.994'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '10.53'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").Value = '2.002.98'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").Value = '1.967'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").Value = '152.21'
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("D28").Value = '17.86'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = '2.030'
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").Value = '112.11'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").Value = '4.878'
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("D32").Value = '0.08691'
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").Value = '3.064'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").Value = '2.746'
$ws.Range("E34").Value = '  +7.41%  '
$ws.Range("D35").Value = '4.454'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").Value = '0.7228'
$ws.Range("E36").Value = '  -4.27%  '
$ws.Range("E37").Value = '  -2.27%  '
$ws.Range("D38").Value = '1.006'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").Value = '1.071'
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").Value = '0.01929'
$ws.Range("D41").Value = '0.05094'
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").Value = '2.868'
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").Value = '6.900'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").Value = '0.5059'
$ws.Range("E44").Value = '  +1.61%  '
$ws.Range("D45").Value = '0.1510'
$ws.Range("E45").Value = '  -5.42%  '
$ws.Range("D46").Value = '8.012'
$ws.Range("E46").Value = '  -3.64%  '
$ws.Range("D47").Value = '1.007'
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").Value = '0.4624'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = '100.95'
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("D50").Value = '9.802'
$ws.Range("E50").Value = '  -3.74%  '
$ws.Range("D51").Value = '1.582'
$ws.Range("E51").Value = '  -2.12%  '
